$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.544.05'
$ws.Range("E2").Value = '  +0.84%  '
$ws.Range("D3").Value = '3.392.19'
$ws.Range("E3").Value = '  +0.00%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.59'
$ws.Range("E5").Value = '  +0.76%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.96'
$ws.Range("E6").Value = '  -0.05%  '
$ws.Range("E8").Value = '  -0.42%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.68'
$ws.Range("E9").Value = '  +0.51%  '
$ws.Range("E10").Value = '  -0.97%  '
$ws.Range("E11").Value = '  -2.17%  '
$ws.Range("D12").Value = '3.967.87'
$ws.Range("E12").Value = '  -0.10%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.126'
$ws.Range("E13").Value = '  +0.04%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.28'
$ws.Range("E14").Value = '  +1.39%  '
$ws.Range("D15").Value = '3.381.44'
$ws.Range("E15").Value = '  -0.17%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000170'
$ws.Range("E16").Value = '  -0.49%  '
$ws.Range("D17").Value = '61.592.89'
$ws.Range("E17").Value = '  +0.83%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.14'
$ws.Range("E18").Value = '  +0.27%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.62'
$ws.Range("E19").Value = '  -1.08%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '9.01'
$ws.Range("E20").Value = '  +1.24%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '390.98'
$ws.Range("E21").Value = '  +2.29%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '74.87'
$ws.Range("E22").Value = '  +0.41%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.549'
$ws.Range("E23").Value = '  -1.26%  '
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("E25").Value = '  +9.35%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000113'
$ws.Range("E26").Value = '  -3.42%  '
$ws.Range("E27").Value = '  +0.13%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.36'
$ws.Range("E28").Value = '  +0.43%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.00'
$ws.Range("E29").Value = '  +0.24%  '
$ws.Range("E30").Value = '  -0.48%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.41'
$ws.Range("E31").Value = '  +0.62%  '
$ws.Range("E32").Value = '  +0.00%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.30'
$ws.Range("E33").Value = '  -0.61%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.91'
$ws.Range("E34").Value = '  -1.00%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '169.36'
$ws.Range("E35").Value = '  +1.46%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.04'
$ws.Range("E36").Value = '  +0.73%  '
$ws.Range("D37").Value = '3.422.38'
$ws.Range("E37").Value = '  -0.08%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.48'
$ws.Range("E38").Value = '  -0.55%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0766'
$ws.Range("E39").Value = '  -0.46%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '25.76'
$ws.Range("E40").Value = '  -4.57%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.779'
$ws.Range("E41").Value = '  -0.26%  '
$ws.Range("E42").Value = '  +0.36%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.65'
$ws.Range("E43").Value = '  -1.35%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.16'
$ws.Range("E44").Value = '  +1.95%  '
$ws.Range("D45").Value = '2.464.79'
$ws.Range("E45").Value = '  -0.40%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.80'
$ws.Range("E46").Value = '  -0.38%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.66'
$ws.Range("E47").Value = '  -1.96%  '
$ws.Range("E48").Value = '  -0.06%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0263'
$ws.Range("E49").Value = '  -1.96%  '
$ws.Range("E50").Value = '  -5.39%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.207'
$ws.Range("E51").Value = '  -1.27%  '
